# Update the date line (title) from "2024-07-11 Thursday" to "2024-07-12 Friday"
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-07-11 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-12 Friday", 2)

# Replace each of the 100 answer cells in the single 20x5 table, in row-major
# order, with their new values (positional replacement -- several old values
# repeat in the table, so this must be done by cell position, not text search).
$newValues = @(
    "6+89=95",
    "41-22=19",
    "91-2=89",
    "90-72=18",
    "72-14=58",
    "45+27=72",
    "60-38=22",
    "85-7=78",
    "91-78=13",
    "75-6=69",
    "39+22=61",
    "39+52=91",
    "66-48=18",
    "19+15=34",
    "51-48=3",
    "27+48=75",
    "37+37=74",
    "7+47=54",
    "26-7=19",
    "42-17=25",
    "61-24=37",
    "81-66=15",
    "55+19=74",
    "38+29=67",
    "83+8=91",
    "36+5=41",
    "90-23=67",
    "61-38=23",
    "64-29=35",
    "38+44=82",
    "33-4=29",
    "7+69=76",
    "8+47=55",
    "80-71=9",
    "47+15=62",
    "83-8=75",
    "34+39=73",
    "93-18=75",
    "46+39=85",
    "74-67=7",
    "14+8=22",
    "65-47=18",
    "72-34=38",
    "48+17=65",
    "82-63=19",
    "69+3=72",
    "66-47=19",
    "45-18=27",
    "92-86=6",
    "81-16=65",
    "74-68=6",
    "61-13=48",
    "53-26=27",
    "9+48=57",
    "93-26=67",
    "36+17=53",
    "3+89=92",
    "82-9=73",
    "92-43=49",
    "97-78=19",
    "77+9=86",
    "56+36=92",
    "14+58=72",
    "39+38=77",
    "12+59=71",
    "19+32=51",
    "63-28=35",
    "34-25=9",
    "13+28=41",
    "90-54=36",
    "92-75=17",
    "82-46=36",
    "80-34=46",
    "92-55=37",
    "18+59=77",
    "13+29=42",
    "28+28=56",
    "16+8=24",
    "37+54=91",
    "7+14=21",
    "13-4=9",
    "9+55=64",
    "26+67=93",
    "81-13=68",
    "54-8=46",
    "57-19=38",
    "83-28=55",
    "24+39=63",
    "86-9=77",
    "96-39=57",
    "34+8=42",
    "63-8=55",
    "30-11=19",
    "3+49=52",
    "63-35=28",
    "95-18=77",
    "71-36=35",
    "69+6=75",
    "86-69=17",
    "17+74=91"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}
